# The post originally stored at row 584 ("「全然分からん」...") was removed.
# Delete that entire row; Excel will automatically shift all subsequent
# rows (585-686) up by one and update the used range / dimension.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(584).Delete()
